$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new data row above the current row 16. This pushes the
# existing rows 16-36 down to 17-37 (dates/values/etc. travel with them),
# and leaves a blank row 16 to be populated with the new record below.
$ws.Range("A16").EntireRow.Insert()

$ws.Range("A16").Value = 11
$ws.Range("B16").Value = "Vega Monumental Concepción"
$ws.Range("C16").Value = "Bíobío"
$ws.Range("D16").Value = 44580
$ws.Range("E16").Value = 8
$ws.Range("F16").Value = 100112031
$ws.Range("G16").Value = "Poroto verde"
$ws.Range("H16").Value = "Magnum"
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 100
$ws.Range("K16").Value = 28000
$ws.Range("L16").Value = 30000
$ws.Range("M16").Value = 29000
$ws.Range("N16").Value = "$/saco 25 kilos"
$ws.Range("O16").Value = "Región Metropolitana"
$ws.Range("P16").Value = 1160
$ws.Range("Q16").Value = 25
$ws.Range("R16").Value = "Hortaliza"
